# Scheduled runner update: refresh market-price derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) across the
# ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 94.833336
$ws.Range("I4").Value = 53.8
$ws.Range("K4").Value = 53.8
$ws.Range("M4").Value = 60.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 24.0625
$ws.Range("I11").Value = 24.0625
$ws.Range("K11").Value = 24.0625
$ws.Range("M11").Value = 115.9375

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1407.3
$ws.Range("I12").Value = 563.6667
$ws.Range("K12").Value = 563.6667
$ws.Range("M12").Value = -393.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 931.3333
$ws.Range("I18").Value = 900
$ws.Range("J18").Value = 994
$ws.Range("K18").Value = 900
$ws.Range("L18").Value = 994
$ws.Range("M18").Value = -616
$ws.Range("N18").Value = -1562

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 37391.375
$ws.Range("J43").Value = 26680
$ws.Range("L43").Value = 26680
$ws.Range("N43").Value = -26818

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1485.2307
$ws.Range("I98").Value = 1268.3636
$ws.Range("J98").Value = 2678
$ws.Range("K98").Value = 1268.3636
$ws.Range("L98").Value = 2678
$ws.Range("M98").Value = 229.6364000000001
$ws.Range("N98").Value = -5674

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1485.2307
$ws.Range("I122").Value = 1268.3636
$ws.Range("J122").Value = 2678
$ws.Range("K122").Value = 3805.0908
$ws.Range("L122").Value = 8034
$ws.Range("M122").Value = -1355.0908
$ws.Range("N122").Value = -12934

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1094.1428
$ws.Range("I125").Value = 833.25
$ws.Range("J125").Value = 1442
$ws.Range("K125").Value = 7499.25
$ws.Range("L125").Value = 12978
$ws.Range("M125").Value = -5039.25
$ws.Range("N125").Value = -17898

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 85780
$ws.Range("J136").Value = 85780
$ws.Range("L136").Value = 85780
$ws.Range("N136").Value = -95980

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 152494.75
$ws.Range("J139").Value = 184989.5
$ws.Range("L139").Value = 184989.5
$ws.Range("N139").Value = -195269.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 547.1923
$ws.Range("I5").Value = 1230.2
$ws.Range("K5").Value = 1230.2
$ws.Range("M5").Value = -1118.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5350.9414
$ws.Range("I32").Value = 4764.148
$ws.Range("K32").Value = 4764.148
$ws.Range("M32").Value = -4477.148

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 34749.168
$ws.Range("J43").Value = 33699
$ws.Range("L43").Value = 33699
$ws.Range("N43").Value = -34325

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 13319.5
$ws.Range("I45").Value = 16225.272
$ws.Range("K45").Value = 16225.272
$ws.Range("M45").Value = -15848.272

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 5074
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 10404.385
$ws.Range("I110").Value = 13562.533
$ws.Range("K110").Value = 13562.533
$ws.Range("M110").Value = -11517.533

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5124.875
$ws.Range("I132").Value = 4666.5
$ws.Range("J132").Value = 6500
$ws.Range("K132").Value = 13999.5
$ws.Range("L132").Value = 19500
$ws.Range("M132").Value = -11469.5
$ws.Range("N132").Value = -24560

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 547.1923
$ws.Range("I4").Value = 1230.2
$ws.Range("K4").Value = 1230.2
$ws.Range("M4").Value = -1115.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 42678.5
$ws.Range("I82").Value = 357
$ws.Range("J82").Value = 85000
$ws.Range("K82").Value = 357
$ws.Range("L82").Value = 85000
$ws.Range("M82").Value = 26
$ws.Range("N82").Value = -85766

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 42678.5
$ws.Range("I85").Value = 357
$ws.Range("J85").Value = 85000
$ws.Range("K85").Value = 357
$ws.Range("L85").Value = 85000
$ws.Range("M85").Value = 969
$ws.Range("N85").Value = -87652

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5047.8945
$ws.Range("I86").Value = 3061.8
$ws.Range("J86").Value = 12495.75
$ws.Range("K86").Value = 3061.8
$ws.Range("L86").Value = 12495.75
$ws.Range("M86").Value = -1938.8
$ws.Range("N86").Value = -14741.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 5047.8945
$ws.Range("I89").Value = 3061.8
$ws.Range("J89").Value = 12495.75
$ws.Range("K89").Value = 15309
$ws.Range("L89").Value = 62478.75
$ws.Range("M89").Value = -9693
$ws.Range("N89").Value = -73710.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 89500
$ws.Range("J132").Value = 89500
$ws.Range("L132").Value = 89500
$ws.Range("N132").Value = -99620

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1572.4706
$ws.Range("I134").Value = 1154.4773
$ws.Range("K134").Value = 3463.4319
$ws.Range("M134").Value = -928.4319

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 259.1
$ws.Range("I22").Value = 256.85715
$ws.Range("J22").Value = 264.33334
$ws.Range("K22").Value = 256.85715
$ws.Range("L22").Value = 264.33334
$ws.Range("M22").Value = 93.14285000000001
$ws.Range("N22").Value = -964.33334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H93").Value = 2407
$ws.Range("I93").Value = 2407
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 2407
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -535
$ws.Range("N93").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H120").Value = 31247.5
$ws.Range("J120").Value = 31247.5
$ws.Range("L120").Value = 31247.5
$ws.Range("N120").Value = -38505.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3163.689
$ws.Range("I132").Value = 2913.476
$ws.Range("J132").Value = 6666.6665
$ws.Range("K132").Value = 8740.428
$ws.Range("L132").Value = 19999.9995
$ws.Range("M132").Value = -6210.428
$ws.Range("N132").Value = -25059.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1459.3871
$ws.Range("I132").Value = 1217.9048
$ws.Range("K132").Value = 10961.1432
$ws.Range("M132").Value = -8431.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 837177.9399999999
$ws.Range("I132").Value = 837177.9399999999
$ws.Range("K132").Value = 2511533.82
$ws.Range("M132").Value = -2509003.82

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1254.4445
$ws.Range("I22").Value = 1378.6
$ws.Range("K22").Value = 1378.6
$ws.Range("M22").Value = -1083.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1254.4445
$ws.Range("I27").Value = 1378.6
$ws.Range("K27").Value = 1378.6
$ws.Range("M27").Value = -1271.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1825
$ws.Range("I46").Value = 1825
$ws.Range("K46").Value = 1825
$ws.Range("M46").Value = -1637

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 53982.715
$ws.Range("I100").Value = 76974.78999999999
$ws.Range("K100").Value = 76974.78999999999
$ws.Range("M100").Value = -76433.78999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4428.7144
$ws.Range("I132").Value = 3182
$ws.Range("K132").Value = 9546
$ws.Range("M132").Value = -7016

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 12375
$ws.Range("I38").Value = 14750
$ws.Range("J38").Value = 10000
$ws.Range("K38").Value = 14750
$ws.Range("L38").Value = 10000
$ws.Range("M38").Value = -14277
$ws.Range("N38").Value = -10946

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2786.3704
$ws.Range("I132").Value = 2911.08
$ws.Range("J132").Value = 1227.5
$ws.Range("K132").Value = 8733.24
$ws.Range("L132").Value = 3682.5
$ws.Range("M132").Value = -6203.24
$ws.Range("N132").Value = -8742.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2241.1875
$ws.Range("I136").Value = 2326.96
$ws.Range("J136").Value = 1934.8572
$ws.Range("K136").Value = 6980.88
$ws.Range("L136").Value = 5804.571599999999
$ws.Range("M136").Value = -4430.88
$ws.Range("N136").Value = -10904.5716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H137").Value = 96489
$ws.Range("I137").Value = 96489
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 96489
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -91389
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H139").Value = 90600
$ws.Range("J139").Value = 90600
$ws.Range("L139").Value = 90600
$ws.Range("N139").Value = -100880

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 113383.336
$ws.Range("J141").Value = 113383.336
$ws.Range("L141").Value = 113383.336
$ws.Range("N141").Value = -123743.336
